$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''63.898.21'
$ws.Range("E2").Value = '  -0.11%  '

$ws.Range("D3").Value = '''2.752.08'
$ws.Range("E3").Value = '  -0.47%  '

$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").Value = '''574.70'
$ws.Range("E5").Value = '  -0.87%  '

$ws.Range("D6").Value = '''157.77'
$ws.Range("E6").Value = '  -0.13%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("E8").Value = '  -1.86%  '

$ws.Range("E9").Value = '  -3.12%  '

$ws.Range("E10").Value = '  +1.29%  '

$ws.Range("D11").Value = '''5.66'
$ws.Range("E11").Value = '  -15.88%  '

$ws.Range("D12").Value = '''0.383'
$ws.Range("E12").Value = '  -3.14%  '

$ws.Range("D13").Value = '''3.238.48'
$ws.Range("E13").Value = '  -0.44%  '

$ws.Range("D14").Value = '''26.57'
$ws.Range("E14").Value = '  -3.63%  '

$ws.Range("D15").Value = '''63.537.71'
$ws.Range("E15").Value = '  -0.63%  '

$ws.Range("E16").Value = '  -2.90%  '

$ws.Range("D17").Value = '''2.756.34'
$ws.Range("E17").Value = '  -0.31%  '

$ws.Range("E18").Value = '  -0.27%  '

$ws.Range("E19").Value = '  -2.63%  '

$ws.Range("E20").Value = '  -2.23%  '

$ws.Range("E21").Value = '  -4.13%  '

$ws.Range("D22").Value = '''1.00'
$ws.Range("E22").Value = '  +0.05%  '

$ws.Range("D23").Value = '''0.534'
$ws.Range("E23").Value = '  -0.68%  '

$ws.Range("D24").Value = '''65.18'
$ws.Range("E24").Value = '  -2.71%  '

$ws.Range("D25").Value = '''0.171'
$ws.Range("E25").Value = '  -0.42%  '

$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  +0.02%  '

$ws.Range("D27").Value = '''8.49'
$ws.Range("E27").Value = '  -1.64%  '

$ws.Range("E28").Value = '  -0.75%  '

$ws.Range("E29").Value = '  -4.85%  '

$ws.Range("D30").Value = '''7.09'
$ws.Range("E30").Value = '  -1.98%  '

$ws.Range("E31").Value = '  -1.51%  '

$ws.Range("D32").Value = '''168.27'
$ws.Range("E32").Value = '  -4.31%  '

$ws.Range("D33").Value = '''20.17'
$ws.Range("E33").Value = '  -2.52%  '

$ws.Range("D34").Value = '''4.92'
$ws.Range("E34").Value = '  -0.14%  '

$ws.Range("E35").Value = '  +0.05%  '

$ws.Range("D36").Value = '''1.46'
$ws.Range("E36").Value = '  +0.18%  '

$ws.Range("D37").Value = '''1.80'
$ws.Range("E37").Value = '  -1.59%  '

$ws.Range("D38").Value = '''0.988'
$ws.Range("E38").Value = '  -2.54%  '

$ws.Range("D39").Value = '''6.16'
$ws.Range("E39").Value = '  +5.90%  '

$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").Value = '''4.15'
$ws.Range("E40").Value = '  -4.10%  '

$ws.Range("B41").Value = 'Bittensor'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D41").Value = '''331.10'
$ws.Range("E41").Value = '  -2.50%  '

$ws.Range("D42").Value = '''38.92'
$ws.Range("E42").Value = '  -1.25%  '

$ws.Range("D43").Value = '''21.49'
$ws.Range("E43").Value = '  -2.16%  '

$ws.Range("E44").Value = '  -2.10%  '

$ws.Range("D45").Value = '''21.56'
$ws.Range("E45").Value = '  -3.81%  '

$ws.Range("E46").Value = '  -2.40%  '

$ws.Range("D47").Value = '''134.94'
$ws.Range("E47").Value = '  -2.05%  '

$ws.Range("D48").Value = '''0.626'
$ws.Range("E48").Value = '  -3.76%  '

$ws.Range("E49").Value = '  -1.00%  '

$ws.Range("E50").Value = '  -0.03%  '

$ws.Range("E51").Value = '  +0.34%  '
